$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Crypto price/volume refresh -- GitHub Actions cron update.
# Force column D cells to Text so numeric-looking price strings (e.g. "1.000",
# "13.70") keep their exact digits instead of being parsed as numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "23.964.98"
$ws.Range("E2").Value = "  +0.56%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.653.87"
$ws.Range("E3").Value = "  +2.20%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.48"
$ws.Range("E5").Value = "  +0.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3893"
$ws.Range("E7").Value = "  -1.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3837"
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("E9").Value = "  +4.15%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.354"
$ws.Range("E10").Value = "  +0.00%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9992"
$ws.Range("E11").Value = "  -0.22%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08452"
$ws.Range("E12").Value = "  +0.01%  "
$ws.Range("E13").Value = "  +0.96%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.122"
$ws.Range("E14").Value = "  +1.06%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.904"
$ws.Range("E15").Value = "  +4.55%  "
$ws.Range("E16").Value = "  +2.72%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.655.03"
$ws.Range("E17").Value = "  +2.58%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.55"
$ws.Range("E18").Value = "  +0.94%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06991"
$ws.Range("E19").Value = "  +0.89%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.76"
$ws.Range("E20").Value = "  -0.95%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.939"
$ws.Range("E21").Value = "  +1.87%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  -0.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.70"
$ws.Range("E23").Value = "  +2.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.000.75"
$ws.Range("E24").Value = "  +0.68%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.473"
$ws.Range("E25").Value = "  +0.56%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.989"
$ws.Range("E26").Value = "  +5.87%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.13"
$ws.Range("E27").Value = "  -0.40%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "151.12"
$ws.Range("E28").Value = "  -3.74%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.445"
$ws.Range("E29").Value = "  +2.65%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "139.23"
$ws.Range("E30").Value = "  -0.89%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.835"
$ws.Range("E31").Value = "  +0.52%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.485"
$ws.Range("E32").Value = "  -0.41%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.836.11"
$ws.Range("E33").Value = "  +2.58%  "
$ws.Range("E34").Value = "  +6.37%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.08113"
$ws.Range("E35").Value = "  -0.09%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02960"
$ws.Range("E36").Value = "  +2.94%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.754"
$ws.Range("E37").Value = "  +2.30%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "10.88"
$ws.Range("E38").Value = "  +5.62%  "
$ws.Range("E39").Value = "  +0.91%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.09149"
$ws.Range("E40").Value = "  +0.15%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.7571"
$ws.Range("E41").Value = "  +0.89%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "13.46"
$ws.Range("E42").Value = "  -1.14%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.429"
$ws.Range("E43").Value = "  +0.20%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.36"
$ws.Range("E44").Value = "  +2.64%  "
$ws.Range("E45").Value = "  +0.58%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.462"
$ws.Range("E46").Value = "  -0.32%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.090"
$ws.Range("E47").Value = "  +0.46%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9996"
$ws.Range("E48").Value = "  -0.07%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.08292"
$ws.Range("E49").Value = "  +0.58%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "134.80"
$ws.Range("E50").Value = "  -0.23%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.216"
